# Added DAP file for today
# Update TOTAL_SUBSTATION_LOAD (col B) and recompute ACTUAL_ENERGY (col D = B - C)
# for hours 10-22 (rows 11-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$loadValues = @{
    11 = 32621
    12 = 33187
    13 = 32997
    14 = 34647
    15 = 35793
    16 = 37060
    17 = 37376
    18 = 37480
    19 = 40256
    20 = 42234
    21 = 40306
    22 = 40825
    23 = 41147
}

foreach ($row in $loadValues.Keys) {
    $ws.Cells.Item($row, 2).Value2 = $loadValues[$row]
    $contestable = $ws.Cells.Item($row, 3).Value2
    $ws.Cells.Item($row, 4).Value2 = $loadValues[$row] - $contestable
}
